# Append: 2025-12-09 06:38 JST
# The scraper re-stamped the "取得日時" (acquisition timestamp) column for
# every existing record row on the "ランサーズ" sheet, moving it from
# 2025-12-09 06:29:13 to 2025-12-09 06:38:11.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$oldTimestamp = "2025-12-09 06:29:13"
$newTimestamp = "2025-12-09 06:38:11"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Value2 -eq $oldTimestamp) {
        $cell.Value = $newTimestamp
    }
}
